$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 6.185
$ws.Range("B10").Value = 5.792
$ws.Range("B12").Value = 5.315
$ws.Range("B18").Value = 5.126
$ws.Range("B37").Value = 8.73
$ws.Range("B55").Value = 4.684
$ws.Range("B68").Value = 5.356
$ws.Range("B77").Value = 6.114
$ws.Range("B78").Value = 7.811
$ws.Range("B81").Value = 6.403999999999999
$ws.Range("B82").Value = 5.659000000000001
